# Apply the "upgrade_to_item_id" column addition to item.xlsx (Sheet1).
# A new column is inserted before the old "expire_time" (F) column, pushing the
# old F (expire_time) -> G and the old G (batch_useable) -> H.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new blank column at F; this shifts F->G and G->H, carrying the
#    existing cell contents/styles with it.
$ws.Columns("F:F").Insert()

# 2) Fill in the new column's header rows (2 and 3 are the schema title rows;
#    4 and 5 stay blank just like they do for the other columns that have no
#    comment/description text).
$ws.Range("F2").Value = "upgrade_to_item_id"
$ws.Range("F3").Value = "int&ref=item.TbItem"

# 3) Fill in the new column's data values (rows 6-15): the id of the item that
#    this row upgrades into.
$ws.Range("F6").Value = 10001
$ws.Range("F7").Value = 10002
$ws.Range("F8").Value = 10003
$ws.Range("F9").Value = 10004
$ws.Range("F10").Value = 10005
$ws.Range("F11").Value = 10006
$ws.Range("F12").Value = 10007
$ws.Range("F13").Value = 10008
$ws.Range("F14").Value = 10009
$ws.Range("F15").Value = 10000

# 4) Adjust column widths: B and E get new explicit widths, and the newly
#    inserted F gets a wide column to fit its long header text. The other
#    columns (A, C, D and the shifted G, H) keep the widths they already had.
$ws.Columns("B:B").ColumnWidth = 15.142857142857142
$ws.Columns("E:E").ColumnWidth = 11.571428571428571
$ws.Columns("F:F").ColumnWidth = 21.714285714285715

# 5) Match the saved selection/active cell.
$ws.Range("D3").Select() | Out-Null
